$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 557.1667
$ws.Range("I11").Value = 557.1667
$ws.Range("K11").Value = 557.1667
$ws.Range("M11").Value = -417.1667
$ws.Range("H12").Value = 199.33333
$ws.Range("I12").Value = 199.33333
$ws.Range("K12").Value = 199.33333
$ws.Range("M12").Value = -29.33332999999999
$ws.Range("H97").Value = 72538.664
$ws.Range("J97").Value = 108333
$ws.Range("L97").Value = 324999
$ws.Range("N97").Value = -325991
$ws.Range("H111").Value = 13065.917
$ws.Range("I111").Value = 1389.4
$ws.Range("J111").Value = 21406.285
$ws.Range("K111").Value = 4168.200000000001
$ws.Range("L111").Value = 64218.855
$ws.Range("M111").Value = -1101.200000000001
$ws.Range("N111").Value = -70352.855
$ws.Range("H112").Value = 38129.93
$ws.Range("I112").Value = 1200
$ws.Range("J112").Value = 44284.918
$ws.Range("K112").Value = 3600
$ws.Range("L112").Value = 132854.754
$ws.Range("M112").Value = -2492
$ws.Range("N112").Value = -135070.754
$ws.Range("H138").Value = 6851913.5
$ws.Range("J138").Value = 9618510
$ws.Range("L138").Value = 28855530
$ws.Range("N138").Value = -28865810

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 23256844
$ws.Range("I2").Value = 31250860
$ws.Range("K2").Value = 31250860
$ws.Range("M2").Value = -31250747
$ws.Range("H97").Value = 1674.4348
$ws.Range("J97").Value = 4249.75
$ws.Range("L97").Value = 4249.75
$ws.Range("N97").Value = -5241.75
$ws.Range("H116").Value = 23256844
$ws.Range("I116").Value = 31250860
$ws.Range("K116").Value = 31250860
$ws.Range("M116").Value = -31248566

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 23256844
$ws.Range("I3").Value = 31250860
$ws.Range("J3").Value = 1529.1818
$ws.Range("K3").Value = 31250860
$ws.Range("L3").Value = 1529.1818
$ws.Range("M3").Value = -31250746
$ws.Range("N3").Value = -1757.1818
$ws.Range("H86").Value = 2222.8572
$ws.Range("I86").Value = 2080.7097
$ws.Range("K86").Value = 2080.7097
$ws.Range("M86").Value = -957.7096999999999
$ws.Range("H89").Value = 2222.8572
$ws.Range("I89").Value = 2080.7097
$ws.Range("K89").Value = 10403.5485
$ws.Range("M89").Value = -4787.548499999999
$ws.Range("H99").Value = 7348.9
$ws.Range("I99").Value = 2299.8
$ws.Range("K99").Value = 2299.8
$ws.Range("M99").Value = -801.8000000000002
$ws.Range("I105").Value = 1911.5
$ws.Range("J105").Value = 1133.75
$ws.Range("K105").Value = 1911.5
$ws.Range("L105").Value = 1133.75
$ws.Range("M105").Value = -164.5
$ws.Range("N105").Value = -4627.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 12887.25
$ws.Range("I23").Value = 2199
$ws.Range("K23").Value = 2199
$ws.Range("M23").Value = -1959
$ws.Range("H27").Value = 12887.25
$ws.Range("I27").Value = 2199
$ws.Range("K27").Value = 2199
$ws.Range("M27").Value = -2007
$ws.Range("H29").Value = 35000
$ws.Range("J29").Value = 35000
$ws.Range("L29").Value = 35000
$ws.Range("N29").Value = -35586
$ws.Range("H31").Value = 58131.777
$ws.Range("I31").Value = 73825.57000000001
$ws.Range("K31").Value = 73825.57000000001
$ws.Range("M31").Value = -73530.57000000001
$ws.Range("H32").Value = 27005
$ws.Range("I32").Value = 27005
$ws.Range("K32").Value = 27005
$ws.Range("M32").Value = -26689
$ws.Range("H34").Value = 58131.777
$ws.Range("I34").Value = 73825.57000000001
$ws.Range("K34").Value = 73825.57000000001
$ws.Range("M34").Value = -73623.57000000001
$ws.Range("H60").Value = 11499.75
$ws.Range("J60").Value = 9000
$ws.Range("L60").Value = 9000
$ws.Range("N60").Value = -10022
$ws.Range("H134").Value = 13312.37
$ws.Range("I134").Value = 3279.7827
$ws.Range("K134").Value = 9839.348100000001
$ws.Range("M134").Value = -7304.348100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 100
$ws.Range("J80").Value = 100
$ws.Range("L80").Value = 300
$ws.Range("N80").Value = -2172
$ws.Range("H81").Value = 5233
$ws.Range("I81").Value = 699
$ws.Range("K81").Value = 2097
$ws.Range("M81").Value = -974
$ws.Range("H83").Value = 100
$ws.Range("J83").Value = 100
$ws.Range("L83").Value = 900
$ws.Range("N83").Value = -10260
$ws.Range("H84").Value = 5233
$ws.Range("I84").Value = 699
$ws.Range("K84").Value = 6291
$ws.Range("M84").Value = -675
$ws.Range("H137").Value = 2993.8572
$ws.Range("I137").Value = 2826.1667
$ws.Range("J137").Value = 4000
$ws.Range("K137").Value = 8478.500100000001
$ws.Range("L137").Value = 12000
$ws.Range("M137").Value = -3378.500100000001
$ws.Range("N137").Value = -22200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 753.63635
$ws.Range("I2").Value = 1150.2142
$ws.Range("J2").Value = 59.625
$ws.Range("K2").Value = 1150.2142
$ws.Range("L2").Value = 59.625
$ws.Range("M2").Value = -1037.2142
$ws.Range("N2").Value = -285.625
$ws.Range("H102").Value = 27779368
$ws.Range("I102").Value = 1670.6177
$ws.Range("K102").Value = 1670.6177
$ws.Range("M102").Value = -48.61770000000001
$ws.Range("H130").Value = 79999
$ws.Range("J130").Value = 79999
$ws.Range("L130").Value = 79999
$ws.Range("N130").Value = -90039
$ws.Range("H132").Value = 2969.4211
$ws.Range("I132").Value = 2745.5
$ws.Range("K132").Value = 8236.5
$ws.Range("M132").Value = -5706.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4050.8286
$ws.Range("I40").Value = 3178.5217
$ws.Range("K40").Value = 3178.5217
$ws.Range("M40").Value = -3042.5217
$ws.Range("H93").Value = 1607.75
$ws.Range("I93").Value = 1599.3636
$ws.Range("J93").Value = 1700
$ws.Range("K93").Value = 1599.3636
$ws.Range("L93").Value = 1700
$ws.Range("M93").Value = -351.3635999999999
$ws.Range("N93").Value = -4196
$ws.Range("H136").Value = 3704.2222
$ws.Range("I136").Value = 3218
$ws.Range("J136").Value = 6500
$ws.Range("K136").Value = 9654
$ws.Range("L136").Value = 19500
$ws.Range("M136").Value = -7104
$ws.Range("N136").Value = -24600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5888.1113
$ws.Range("I81").Value = 8002.2144
$ws.Range("K81").Value = 16004.4288
$ws.Range("M81").Value = -14943.4288
$ws.Range("H84").Value = 5888.1113
$ws.Range("I84").Value = 8002.2144
$ws.Range("K84").Value = 80022.144
$ws.Range("M84").Value = -74718.144
$ws.Range("H122").Value = 2054.889
$ws.Range("I122").Value = 1914.6154
$ws.Range("J122").Value = 2419.6
$ws.Range("K122").Value = 5743.8462
$ws.Range("L122").Value = 7258.799999999999
$ws.Range("M122").Value = -3293.8462
$ws.Range("N122").Value = -12158.8
